{"js": "// The \"Flujos Alternativos\" step-number tables had a couple of step\n// columns that were mis-numbered (continuing a shared numbering restart).\n// This renumbers the trailing steps in the affected tables so the\n// sequence reads on without repeats/gaps:\n//   - \"El correo ya existe en el sistema.\" table: steps 4,5,6,7 -> 5,6,7,8\n//     (the row that already read \"4\" stays \"4\"; the *next* row, which\n//     incorrectly repeated \"4\", becomes \"5\", and every following row\n//     shifts up by one).\n//   - \"Campos incompletos:\" table: steps 2,3,4,5,6 -> 4,5,6,7,8\n//     (after the row that already read \"3\", the flow had restarted at\n//     \"2\" and needs to continue at \"4\").\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Each entry: [tableIndex (0-based), rowIndex (0-based), expectedOldText, newText]\nconst renumbers = [\n  [1, 5, \"4\", \"5\"],\n  [1, 6, \"5\", \"6\"],\n  [1, 7, \"6\", \"7\"],\n  [1, 8, \"7\", \"8\"],\n  [3, 4, \"2\", \"4\"],\n  [3, 5, \"3\", \"5\"],\n  [3, 6, \"4\", \"6\"],\n  [3, 7, \"5\", \"7\"],\n  [3, 8, \"6\", \"8\"],\n];\n\nfor (const [tableIndex, rowIndex, oldText, newText] of renumbers) {\n  const table = tables.items[tableIndex];\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  const row = rows.items[rowIndex];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  const cell = cells.items[0];\n  cell.load(\"value\");\n  await context.sync();\n\n  if (cell.value.trim() === oldText) {\n    cell.body.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The \"Flujos Alternativos\" step-number tables had a couple of step\n# columns that were mis-numbered (continuing a shared numbering restart).\n# This renumbers the trailing steps in the affected tables so the\n# sequence reads on without repeats/gaps:\n#   - \"El correo ya existe en el sistema.\" table (2nd table): steps\n#     4,5,6,7 -> 5,6,7,8 (the row that already read \"4\" stays \"4\"; the\n#     *next* row, which incorrectly repeated \"4\", becomes \"5\", and every\n#     following row shifts up by one).\n#   - \"Campos incompletos:\" table (4th table): steps 2,3,4,5,6 -> 4,5,6,7,8\n#     (after the row that already read \"3\", the flow had restarted at\n#     \"2\" and needs to continue at \"4\").\n\n$d = $word.ActiveDocument\n\n# Each entry: table index (1-based), row index (1-based), expected old text, new text\n$renumbers = @(\n    @(2, 6, \"4\", \"5\"),\n    @(2, 7, \"5\", \"6\"),\n    @(2, 8, \"6\", \"7\"),\n    @(2, 9, \"7\", \"8\"),\n    @(4, 5, \"2\", \"4\"),\n    @(4, 6, \"3\", \"5\"),\n    @(4, 7, \"4\", \"6\"),\n    @(4, 8, \"5\", \"7\"),\n    @(4, 9, \"6\", \"8\")\n)\n\nforeach ($item in $renumbers) {\n    $tableIndex = $item[0]\n    $rowIndex = $item[1]\n    $oldText = $item[2]\n    $newText = $item[3]\n\n    $table = $d.Tables.Item($tableIndex)\n    $cell = $table.Cell($rowIndex, 1)\n    # Cell.Range.Text carries a trailing paragraph mark (chr 13) and\n    # cell mark (chr 7); strip those before comparing the visible text.\n    $current = $cell.Range.Text.TrimEnd([char]7, [char]13)\n\n    if ($current -eq $oldText) {\n        $cell.Range.Text = $newText\n    }\n}\n"}
